$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.305.85"
$ws.Range("E2").Value = "  -2.12%  "
$ws.Range("D3").Value = "1.825.77"
$ws.Range("E3").Value = "  -1.88%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -1.34%  "
$ws.Range("D5").Value = "'314.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.09%  "
$ws.Range("D6").Value = "'1.004"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.17%  "
$ws.Range("D7").Value = "'0.4276"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.06%  "
$ws.Range("D8").Value = "'0.3693"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.88%  "
$ws.Range("D9").Value = "'0.07257"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.33%  "
$ws.Range("E10").Value = "  -2.59%  "
$ws.Range("D11").Value = "'21.08"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.12%  "
$ws.Range("D12").Value = "1.807.76"
$ws.Range("E12").Value = "  -2.91%  "
$ws.Range("D13").Value = "'6.715"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.04%  "
$ws.Range("D14").Value = "'0.07105"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.43%  "
$ws.Range("D15").Value = "'5.319"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.48%  "
$ws.Range("D16").Value = "'88.76"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("E17").Value = "  -1.32%  "
$ws.Range("D18").Value = "'0.000008873"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.10%  "
$ws.Range("D19").Value = "'1.004"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.16%  "
$ws.Range("E20").Value = "  -3.07%  "
$ws.Range("D21").Value = "27.326.61"
$ws.Range("E21").Value = "  -2.07%  "
$ws.Range("D22").Value = "'5.150"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.64%  "
$ws.Range("D23").Value = "'10.88"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.06%  "
$ws.Range("D24").Value = "2.048.85"
$ws.Range("E24").Value = "  -2.20%  "
$ws.Range("D25").Value = "'2.007"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.29%  "
$ws.Range("E26").Value = "  -2.48%  "
$ws.Range("D27").Value = "'18.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.77%  "
$ws.Range("D28").Value = "'2.149"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.90%  "
$ws.Range("D29").Value = "'5.256"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.32%  "
$ws.Range("D30").Value = "'116.83"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.74%  "
$ws.Range("D31").Value = "'0.08888"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.17%  "
$ws.Range("D32").Value = "'1.204"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.37%  "
$ws.Range("D33").Value = "'0.7592"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.59%  "
$ws.Range("D34").Value = "'4.459"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.11%  "
$ws.Range("D35").Value = "'2.840"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.30%  "
$ws.Range("E36").Value = "  -1.23%  "
$ws.Range("D37").Value = "'1.115"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.37%  "
$ws.Range("D38").Value = "'0.01985"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.51%  "
$ws.Range("D39").Value = "'0.05285"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.87%  "
$ws.Range("D40").Value = "'7.171"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.14%  "
$ws.Range("D41").Value = "'2.874"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.36%  "
$ws.Range("E42").Value = "  +0.73%  "
$ws.Range("D43").Value = "'0.5056"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.07%  "
$ws.Range("D44").Value = "'8.698"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.05%  "
$ws.Range("D45").Value = "'10.68"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.01%  "
$ws.Range("D46").Value = "'107.78"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.85%  "
$ws.Range("D47").Value = "'0.4761"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.06%  "
$ws.Range("E48").Value = "  -1.26%  "
$ws.Range("D49").Value = "'0.06371"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.41%  "
$ws.Range("D50").Value = "'1.667"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.29%  "
$ws.Range("D51").Value = "'1.841"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.32%  "
